$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 169.915657
$ws.Range("H2").Value = 509.746971
$ws.Range("I2").Value = 0.4441184931734509
$ws.Range("J2").Value = 0.4441184931734509
$ws.Range("M2").Value = 0.007957000000000001
$ws.Range("N2").Value = 0.023871
$ws.Range("O2").Value = 0.0002448939493579708
$ws.Range("P2").Value = 0.0002448939493579708
$ws.Range("Q2").Value = 1.352018882749
$ws.Range("R2").Value = 12.168169944741
$ws.Range("S2").Value = 0.0001087619317761574
$ws.Range("T2").Value = 0.0001087619317761574
$ws.Range("G3").Value = 169.915657
$ws.Range("H3").Value = 509.746971
$ws.Range("I3").Value = 0.4441184931734509
$ws.Range("J3").Value = 0.4441184931734509
$ws.Range("O3").Value = 0.003249135679578298
$ws.Range("P3").Value = 0.003249135679578299
$ws.Range("Q3").Value = 17.93793927093767
$ws.Range("R3").Value = 161.441453438439
$ws.Range("S3").Value = 0.00144300124213041
$ws.Range("T3").Value = 0.001443001242130411
$ws.Range("G4").Value = 169.915657
$ws.Range("H4").Value = 509.746971
$ws.Range("I4").Value = 0.4441184931734509
$ws.Range("J4").Value = 0.4441184931734509
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02449766666666667
$ws.Range("N4").Value = 0.073493
$ws.Range("O4").Value = 0.0007539688752111494
$ws.Range("P4").Value = 0.0007539688752111494
$ws.Range("Q4").Value = 4.162537126633667
$ws.Range("R4").Value = 37.46283413970301
$ws.Range("S4").Value = 0.0003348515207584573
$ws.Range("T4").Value = 0.0003348515207584573
$ws.Range("G5").Value = 169.915657
$ws.Range("H5").Value = 509.746971
$ws.Range("I5").Value = 0.4441184931734509
$ws.Range("J5").Value = 0.4441184931734509
$ws.Range("M5").Value = 32.353591
$ws.Range("N5").Value = 97.060773
$ws.Range("O5").Value = 0.9957520014958525
$ws.Range("P5").Value = 0.9957520014958525
$ws.Range("Q5").Value = 5497.381671074288
$ws.Range("R5").Value = 49476.43503966858
$ws.Range("S5").Value = 0.4422318784787858
$ws.Range("T5").Value = 0.4422318784787859
$ws.Range("I6").Value = 0.1787346690539575
$ws.Range("J6").Value = 0.1787346690539575
$ws.Range("M6").Value = 0.007957000000000001
$ws.Range("N6").Value = 0.023871
$ws.Range("O6").Value = 0.0002448939493579708
$ws.Range("P6").Value = 0.0002448939493579708
$ws.Range("Q6").Value = 0.5441175075510001
$ws.Range("R6").Value = 4.897057567959
$ws.Range("S6").Value = 0.00004377103899181353
$ws.Range("T6").Value = 0.00004377103899181353
$ws.Range("I7").Value = 0.1787346690539575
$ws.Range("J7").Value = 0.1787346690539575
$ws.Range("O7").Value = 0.003249135679578298
$ws.Range("P7").Value = 0.003249135679578299
$ws.Range("S7").Value = 0.0005807331904008325
$ws.Range("T7").Value = 0.0005807331904008326
$ws.Range("I8").Value = 0.1787346690539575
$ws.Range("J8").Value = 0.1787346690539575
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.02449766666666667
$ws.Range("N8").Value = 0.073493
$ws.Range("O8").Value = 0.0007539688752111494
$ws.Range("P8").Value = 0.0007539688752111494
$ws.Range("Q8").Value = 1.675205394933
$ws.Range("R8").Value = 15.076848554397
$ws.Range("S8").Value = 0.0001347603773878494
$ws.Range("T8").Value = 0.0001347603773878494
$ws.Range("I9").Value = 0.1787346690539575
$ws.Range("J9").Value = 0.1787346690539575
$ws.Range("M9").Value = 32.353591
$ws.Range("N9").Value = 97.060773
$ws.Range("O9").Value = 0.9957520014958525
$ws.Range("P9").Value = 0.9957520014958525
$ws.Range("Q9").Value = 2212.411121684613
$ws.Range("R9").Value = 19911.70009516152
$ws.Range("S9").Value = 0.177975404447177
$ws.Range("T9").Value = 0.177975404447177
$ws.Range("G10").Value = 53.27463399999999
$ws.Range("H10").Value = 159.823902
$ws.Range("I10").Value = 0.1392470275793777
$ws.Range("J10").Value = 0.1392470275793778
$ws.Range("M10").Value = 0.007957000000000001
$ws.Range("N10").Value = 0.023871
$ws.Range("O10").Value = 0.0002448939493579708
$ws.Range("P10").Value = 0.0002448939493579708
$ws.Range("Q10").Value = 0.423906262738
$ws.Range("R10").Value = 3.815156364641999
$ws.Range("S10").Value = 0.00003410075452027209
$ws.Range("T10").Value = 0.0000341007545202721
$ws.Range("G11").Value = 53.27463399999999
$ws.Range("H11").Value = 159.823902
$ws.Range("I11").Value = 0.1392470275793777
$ws.Range("J11").Value = 0.1392470275793778
$ws.Range("O11").Value = 0.003249135679578298
$ws.Range("P11").Value = 0.003249135679578299
$ws.Range("Q11").Value = 5.624185353168667
$ws.Range("R11").Value = 50.617668178518
$ws.Range("S11").Value = 0.0004524324855833796
$ws.Range("T11").Value = 0.0004524324855833797
$ws.Range("G12").Value = 53.27463399999999
$ws.Range("H12").Value = 159.823902
$ws.Range("I12").Value = 0.1392470275793777
$ws.Range("J12").Value = 0.1392470275793778
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.02449766666666667
$ws.Range("N12").Value = 0.073493
$ws.Range("O12").Value = 0.0007539688752111494
$ws.Range("P12").Value = 0.0007539688752111494
$ws.Range("Q12").Value = 1.305104225520666
$ws.Range("R12").Value = 11.745938029686
$ws.Range("S12").Value = 0.0001049879247605193
$ws.Range("T12").Value = 0.0001049879247605193
$ws.Range("G13").Value = 53.27463399999999
$ws.Range("H13").Value = 159.823902
$ws.Range("I13").Value = 0.1392470275793777
$ws.Range("J13").Value = 0.1392470275793778
$ws.Range("M13").Value = 32.353591
$ws.Range("N13").Value = 97.060773
$ws.Range("O13").Value = 0.9957520014958525
$ws.Range("P13").Value = 0.9957520014958525
$ws.Range("Q13").Value = 1723.625719110694
$ws.Range("R13").Value = 15512.63147199624
$ws.Range("S13").Value = 0.1386555064145135
$ws.Range("T13").Value = 0.1386555064145136
$ws.Range("G14").Value = 91.01828266666666
$ws.Range("H14").Value = 273.054848
$ws.Range("I14").Value = 0.2378998101932138
$ws.Range("J14").Value = 0.2378998101932138
$ws.Range("M14").Value = 0.007957000000000001
$ws.Range("N14").Value = 0.023871
$ws.Range("O14").Value = 0.0002448939493579708
$ws.Range("P14").Value = 0.0002448939493579708
$ws.Range("Q14").Value = 0.7242324751786667
$ws.Range("R14").Value = 6.518092276608
$ws.Range("S14").Value = 0.00005826022406972775
$ws.Range("T14").Value = 0.00005826022406972776
$ws.Range("G15").Value = 91.01828266666666
$ws.Range("H15").Value = 273.054848
$ws.Range("I15").Value = 0.2378998101932138
$ws.Range("J15").Value = 0.2378998101932138
$ws.Range("O15").Value = 0.003249135679578298
$ws.Range("P15").Value = 0.003249135679578299
$ws.Range("Q15").Value = 9.608769761692445
$ws.Range("R15").Value = 86.478927855232
$ws.Range("S15").Value = 0.0007729687614636759
$ws.Range("T15").Value = 0.000772968761463676
$ws.Range("G16").Value = 91.01828266666666
$ws.Range("H16").Value = 273.054848
$ws.Range("I16").Value = 0.2378998101932138
$ws.Range("J16").Value = 0.2378998101932138
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.02449766666666667
$ws.Range("N16").Value = 0.073493
$ws.Range("O16").Value = 0.0007539688752111494
$ws.Range("P16").Value = 0.0007539688752111494
$ws.Range("Q16").Value = 2.229735549340444
$ws.Range("R16").Value = 20.067619944064
$ws.Range("S16").Value = 0.0001793690523043233
$ws.Range("T16").Value = 0.0001793690523043234
$ws.Range("G17").Value = 91.01828266666666
$ws.Range("H17").Value = 273.054848
$ws.Range("I17").Value = 0.2378998101932138
$ws.Range("J17").Value = 0.2378998101932138
$ws.Range("M17").Value = 32.353591
$ws.Range("N17").Value = 97.060773
$ws.Range("O17").Value = 0.9957520014958525
$ws.Range("P17").Value = 0.9957520014958525
$ws.Range("Q17").Value = 2944.768290919723
$ws.Range("R17").Value = 26502.9146182775
$ws.Range("S17").Value = 0.2368892121553761
$ws.Range("T17").Value = 0.2368892121553761
